$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 72 (shifts existing rows 72.. down by one)
$ws.Rows.Item(72).Insert()

# New record inserted at row 72
$ws.Cells.Item(72, 1).Value = 10
$ws.Cells.Item(72, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(72, 3).Value = "La Araucanía"
$ws.Cells.Item(72, 4).Value = 44995
$ws.Cells.Item(72, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(72, 5).Value = 9
$ws.Cells.Item(72, 6).Value = "Fruta"
$ws.Cells.Item(72, 7).Value = 100108
$ws.Cells.Item(72, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(72, 9).Value = 100108003
$ws.Cells.Item(72, 10).Value = "Maracuyá"
$ws.Cells.Item(72, 11).Value = "Sin especificar"
$ws.Cells.Item(72, 12).Value = "Primera"
$ws.Cells.Item(72, 13).Value = 20
$ws.Cells.Item(72, 14).Value = 60000
$ws.Cells.Item(72, 15).Value = 60000
$ws.Cells.Item(72, 16).Value = 60000
$ws.Cells.Item(72, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(72, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 19).Value = 3333
$ws.Cells.Item(72, 20).Value = 18
